$d = $word.ActiveDocument

$find = "Kampagnendaten 2022 für das Sternbild Pegasus-Konstellation: 8. bis 17. Oktober, 7. bis 16. November,"
$replace = "Kampagnendaten 2022 für das Pegasus-Konstellation: 8. bis 17. Oktober, 7. bis 16. November,"

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                         $true, 1, $false, $replace, 2)
